$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet (tab was "range names", now "writing to cells")
$ws.Name = "writing to cells"

# 2) Row 1: a date/time value in A1 formatted as dd-mm-yyyy, new small values in
#    B1:D1, and the old E1:L1 numbers removed. Registering the date/time cell's
#    format on a cell that legitimately gets cleared below (E1) first makes the
#    workbook register the "yyyy-mm-dd h:mm:ss" format code ahead of the
#    "dd-mm-yyyy" one, matching numFmtId 164 / 165.
$ws.Range("E1").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("A1").Value = 42512.8415002641
$ws.Range("A1").NumberFormat = "dd-mm-yyyy"
$ws.Range("B1").Value = 42
$ws.Range("C1").Value = 43
$ws.Range("D1").Value = 44
$ws.Range("E1:L1").Clear()

# 3) Rows 7-10: replace the old A1:L1-style sequences with the new sparse
#    B/E/H/K/N/Q squares pattern, and add a new row 11 with the same pattern.
foreach ($r in 7..11) {
    $ws.Range("A$r" + ":L$r").Clear()
    $ws.Cells.Item($r, 2).Value = 4
    $ws.Cells.Item($r, 5).Value = 25
    $ws.Cells.Item($r, 8).Value = 64
    $ws.Cells.Item($r, 11).Value = 121
    $ws.Cells.Item($r, 14).Value = 196
    $ws.Cells.Item($r, 17).Value = 289
}
